$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record (Jengibre, 2021-12-16) needs to be inserted
# as row 34, pushing the existing rows 34-59 down to 35-60.
$ws.Rows.Item(34).Insert()

# Populate the newly inserted row 34 with the new record's data. Columns
# that are identical for every Jengibre record at this market are copied
# from the surrounding rows; only Fecha/Volumen/Precios differ.
$ws.Cells.Item(34, 1).Value() = 6
$ws.Cells.Item(34, 2).Value() = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(34, 3).Value() = "Metropolitana"
$ws.Cells.Item(34, 4).Value() = 44546
$ws.Cells.Item(34, 5).Value() = 13
$ws.Cells.Item(34, 6).Value() = 100114007
$ws.Cells.Item(34, 7).Value() = "Jengibre"
$ws.Cells.Item(34, 8).Value() = "Sin especificar"
$ws.Cells.Item(34, 9).Value() = "Primera"
$ws.Cells.Item(34, 10).Value() = 200
$ws.Cells.Item(34, 11).Value() = 13000
$ws.Cells.Item(34, 12).Value() = 15000
$ws.Cells.Item(34, 13).Value() = 14200
$ws.Cells.Item(34, 14).Value() = "$/caja 13 kilos"
$ws.Cells.Item(34, 15).Value() = "Perú"
$ws.Cells.Item(34, 16).Value() = 1092
$ws.Cells.Item(34, 17).Value() = 13
$ws.Cells.Item(34, 18).Value() = "Hortaliza"
